$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.6545652718822623
$ws.Range("C2").Value = 1.626987699542094
$ws.Range("D2").Value = 0.1496068669990043
$ws.Range("E2").Value = 0.5333859586016987
$ws.Range("G2").Value = 2.964545797025059

$ws.Range("B3").Value = 0.2881169905109251
$ws.Range("C3").Value = 0.3048912486333797
$ws.Range("D3").Value = 189.6080260415259
$ws.Range("E3").Value = 0.5333859586016987
$ws.Range("G3").Value = 190.7344202392719

$ws.Range("B4").Value = 0.1169995834814548
$ws.Range("C4").Value = 0.3048912486333797
$ws.Range("D4").Value = 0.7210945179870265
$ws.Range("E4").Value = 0.5333859586016987
$ws.Range("G4").Value = 1.67637130870356

$ws.Range("B5").Value = 3.272327238179451
$ws.Range("C5").Value = 1.626987699542094
$ws.Range("D5").Value = 0.7210945179870265
$ws.Range("E5").Value = 0.5333859586016987
$ws.Range("G5").Value = 6.15379541431027

$ws.Range("B6").Value = 1.445647641019636
$ws.Range("C6").Value = 1.626987699542094
$ws.Range("D6").Value = 3.223369029078222
$ws.Range("E6").Value = 0.5333859586016987
$ws.Range("G6").Value = 6.82939032824165

$ws.Range("B7").Value = 3.272327238179451
$ws.Range("C7").Value = 1.626987699542094
$ws.Range("D7").Value = 0.7210945179870265
$ws.Range("E7").Value = 0.5333859586016987
$ws.Range("G7").Value = 6.15379541431027

$ws.Range("B8").Value = 0.1169995834814548
$ws.Range("C8").Value = 0.3048912486333797
$ws.Range("D8").Value = 0.7210945179870265
$ws.Range("E8").Value = 0.5333859586016987
$ws.Range("G8").Value = 1.67637130870356

$ws.Range("B9").Value = 3.272327238179451
$ws.Range("C9").Value = 1.626987699542094
$ws.Range("D9").Value = 3.223369029078222
$ws.Range("E9").Value = 0.5333859586016987
$ws.Range("G9").Value = 8.656069925401464

$ws.Range("B10").Value = 0.1169995834814548
$ws.Range("C10").Value = 0.3048912486333797
$ws.Range("D10").Value = 0.1496068669990043
$ws.Range("E10").Value = 0.5333859586016987
$ws.Range("G10").Value = 1.104883657715537

$ws.Range("B11").Value = 3.272327238179451
$ws.Range("C11").Value = 1.626987699542094
$ws.Range("D11").Value = 0.1496068669990043
$ws.Range("E11").Value = 0.5333859586016987
$ws.Range("G11").Value = 5.582307763322248

$ws.Range("B12").Value = 3.272327238179451
$ws.Range("C12").Value = 1.626987699542094
$ws.Range("D12").Value = 0.7210945179870265
$ws.Range("E12").Value = 0.5333859586016987
$ws.Range("G12").Value = 6.15379541431027

$ws.Range("B13").Value = 0.6545652718822623
$ws.Range("C13").Value = 1.626987699542094
$ws.Range("D13").Value = 3.223369029078222
$ws.Range("E13").Value = 0.5333859586016987
$ws.Range("G13").Value = 6.038307959104277

$ws.Range("B14").Value = 3.272327238179451
$ws.Range("C14").Value = 1.626987699542094
$ws.Range("D14").Value = 0.1496068669990043
$ws.Range("E14").Value = 0.5333859586016987
$ws.Range("G14").Value = 5.582307763322248

$ws.Range("B15").Value = 0.6545652718822623
$ws.Range("C15").Value = 1.626987699542094
$ws.Range("D15").Value = 0.7210945179870265
$ws.Range("E15").Value = 0.5333859586016987
$ws.Range("G15").Value = 3.536033448013082

$ws.Range("B16").Value = 0.6545652718822623
$ws.Range("C16").Value = 0.3048912486333797
$ws.Range("D16").Value = 0.1496068669990043
$ws.Range("E16").Value = 0.5333859586016987
$ws.Range("G16").Value = 1.642449346116345

$ws.Range("B17").Value = 1.445647641019636
$ws.Range("C17").Value = 1.626987699542094
$ws.Range("D17").Value = 0.7210945179870265
$ws.Range("E17").Value = 0.5333859586016987
$ws.Range("G17").Value = 4.327115817150455
